# The workbook originally has 7 sheets:
#   1: 5b1fa97a-a26d-3695-a   (note "18." - Contract assets)
#   2: 6179dbd7-5cac-33ab-9
#   3: 29c71f19-706c-3b96-9
#   4: 1c24c1ad-8ed2-3d96-8
#   5: 1aa9caad-49d7-3aa1-a
#   6: 15fa5e32-2817-3795-a
#   7: ad82a938-84a8-34e3-9
#
# The fix inserts a new note "17. Trade and other receivables" sheet in front
# of the existing "18." sheet (which keeps its original content, but is
# pushed down one position), and renumbers/prefixes every tab name.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

# Step 1: duplicate the current sheet 1 ("18." note) and drop the copy right
# after sheet 4, so it becomes the new sheet 5 -- preserving its original
# content untouched.
$ws1.Copy($null, $ws4)

# Keep the original sheet 1 as the active/selected tab (matches the source
# workbook, where the first sheet was the selected one).
$ws1.Activate()

# Step 2: turn the original sheet 1 into the new "17." note (trade and other
# receivables) with the corrected figures.
$ws1.Range("A2").Value = "'17."
$ws1.Range("A3").Value = "Int thousands of AUD"
$ws1.Range("B3").Value = "2022"
$ws1.Range("C3").Value = "2021"

$ws1.Range("A4").Value = "Amount expected to be settled within 12 months"
$ws1.Range("B4").Value = "'8810"
$ws1.Range("C4").Value = "'6677"

$ws1.Range("A5").Value = "Amount expected to be settled after more than 12 months"

$ws1.Range("B6").Value = "'8810"
$ws1.Range("C6").Value = "'6677"

$ws1.Range("A7").Value = "Reconciliation at the beginning and end of the current and previous financial year are set out below"

$ws1.Range("A8").Value = "Opening Balance"
$ws1.Range("B8").Value = "'6677"
$ws1.Range("C8").Value = "'6255"

$ws1.Range("A9").Value = "Transfer to trade receivables included in opening balance"
$ws1.Range("B9").Value = "(6,677)"
$ws1.Range("C9").Value = "(6,255)"

$ws1.Range("A10").Value = "Additions"
$ws1.Range("B10").Value = "'8810"
$ws1.Range("C10").Value = "'6677"

$ws1.Range("A11").Value = "Closing balance"
$ws1.Range("B11").Value = "'8810"
$ws1.Range("C11").Value = "'6677"

# Step 3: rename every tab to match the new naming scheme.
$wb.Worksheets.Item(1).Name = "17__5b1fa97a-a26d-36"
$wb.Worksheets.Item(2).Name = "15__6179dbd7-5cac-33"
$wb.Worksheets.Item(3).Name = "15__29c71f19-706c-3b"
$wb.Worksheets.Item(4).Name = "10__1c24c1ad-8ed2-3d"
$wb.Worksheets.Item(5).Name = "18__5b1fa97a-a26d-36"
$wb.Worksheets.Item(6).Name = "11__1aa9caad-49d7-3a"
$wb.Worksheets.Item(7).Name = "12__15fa5e32-2817-37"
$wb.Worksheets.Item(8).Name = "12__ad82a938-84a8-34"
